# Applies the "Add files via upload" revision to Makenson Noel's resume.
# The underlying change is a single real content edit (adding
# " Resort & Spa" after "Marriott Harbor Beach"); every other hunk in the
# source diff only reflects Word re-flowing / merging previously-split
# <w:r> runs that already carried identical text - no visible wording
# changes there. We reproduce both: the one real text insertion, and the
# run-merges (which happen naturally once Find/Replace rewrites the text).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Education: Bachelor's "Related coursework" line (run merge only) ---
Replace-Text "Related coursework: Data Structures, Stochastic Models, Database Structures, Computer Operating Systems, Artificial Intelligence, Mobile Apps, Principles of Software Engineering, Calculus II, Physics II, Foundations of Cybersecurity, Deep Learning, Internet of things, Engineering Design I." `
             "Related coursework: Data Structures, Stochastic Models, Database Structures, Computer Operating Systems, Artificial Intelligence, Mobile Apps, Principles of Software Engineering, Calculus II, Physics II, Foundations of Cybersecurity, Deep Learning, Internet of things, Engineering Design I."

# --- Education: Associate's "Major: Computer Science" (run merge only) ---
# "Major: Computer Science" also appears (already merged) under the
# Bachelor's degree, so anchor the search after "Associate of Arts" to hit
# only the Associate's-degree occurrence that is still split into runs.
$anchor = $d.Content
$null = $anchor.Find.Execute("Associate of Arts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scoped = $d.Range($anchor.End, $d.Content.End)
$null = $scoped.Find.Execute("Major: Comput" + "er Science", $true, $false, $false, $false, $false, $true, 1, $false, "Major: Computer Science", 2)

# --- Experience: Food Runner | Marriott Harbor Beach (+ " Resort & Spa") ---
$rng = $d.Content
$null = $rng.Find.Execute("Marriott Harbor Beach", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.InsertAfter(" Resort & Spa")

# --- Experience: "Worked in hospitality..." (run merge only) ---
Replace-Text "Worked in hospitality at Marriott for 4 Years while working towards a Bachelor of Science in Computer Science Degree. This experience improved time management, teamwork, and communication skills." `
             "Worked in hospitality at Marriott for 4 Years while working towards a Bachelor of Science in Computer Science Degree. This experience improved time management, teamwork, and communication skills."

# --- Projects: COVID-19 App description (run merge only) ---
Replace-Text "A Coronavirus 19 (COVID-19) application was created to inform users on COVID-19. The application was developed using Cordova, Html, JavaScript, and Bootstrap. " `
             "A Coronavirus 19 (COVID-19) application was created to inform users on COVID-19. The application was developed using Cordova, Html, JavaScript, and Bootstrap. "

# --- Projects: "Coding is Fun | May 2019 - June 2019" (run merge only) ---
# Include the trailing space so the whole second run (" 2019 - June 2019 ")
# is absorbed into the merge, matching the target's single trailing-space run.
Replace-Text "Fun | May 2019 – June 2019 " "Fun | May 2019 – June 2019 "

# --- Organizations: "Inducted in May 2017..." (run merge only) ---
Replace-Text "Inducted in May 2017, received the National Society of Leadership and Success Certificate." `
             "Inducted in May 2017, received the National Society of Leadership and Success Certificate."

# --- Technical Skills: "DaTabase Management" heading (run merge only) ---
Replace-Text ("Da" + "T" + "abase Management") "DaTabase Management"

# --- Technical Skills: "MySQL, Oracle Database, phpMyAdmin. " (run merge only) ---
Replace-Text "MySQL, Oracle Database, phpMyAdmin. " "MySQL, Oracle Database, phpMyAdmin. "

# --- "Volunteering Experience" heading (run merge only) ---
Replace-Text ("V" + "olunteering Experience") "Volunteering Experience"

# --- Volunteering: "Give Kids The World Village | December 2020" (run merge + drop proofErr) ---
# Include the trailing space in the match so the found range spans the
# <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/> markers
# that used to straddle "2020" - Replace collapses them away along with
# the run split.
Replace-Text "Give Kids The World Village | December 2020 " "Give Kids The World Village | December 2020 "

# --- Volunteering: "Volunteered for Give Kids the World Village..." (run merge only) ---
Replace-Text "Volunteered for “Give Kids the World Village” in Kissimmee, Florida, it is a nonprofit resort for children with critical illnesses and their families." `
             "Volunteered for “Give Kids the World Village” in Kissimmee, Florida, it is a nonprofit resort for children with critical illnesses and their families."

# --- Volunteering: "Was tasked with taking photos..." (run merge only) ---
Replace-Text "Was tasked with taking photos of visitors, encouraging them to press the wish button on the wish tree, sanitize their hands and practice social distancing." `
             "Was tasked with taking photos of visitors, encouraging them to press the wish button on the wish tree, sanitize their hands and practice social distancing."
